$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new row of data describing the "Switch spell_pos" message
$ws.Range("B7").Value = 5
$ws.Range("C7").Value = "Switch spell_pos"
$ws.Range("D7").Value = "(weapon/idx)(weapon/idx)"

# Widen column D to fit the new content
# (target stored width is 22.33203125 chars; the COM width model here
# quantizes to 1/6-character steps, so 21.5 is the closest achievable
# ColumnWidth input and lands on 22.33333... in the saved XML)
$ws.Columns.Item(4).ColumnWidth = 21.5

# Update the selected cell to reflect where the user ended up after editing
$ws.Range("D11").Select()
